# Adding the changes we made on may 9th
# Insert 10 new data rows right after the header row (row 1), pushing the
# existing data down by 10 rows, then populate those new rows with the
# new data points.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 10 rows at row 2 (the existing row 2..21 data will shift to 12..31)
$ws.Rows("2:11").Insert()
# The insert copies formatting down from the header row; clear it so the new
# data rows stay unformatted like the rest of the numeric rows.
$ws.Range("A2:C11").ClearFormats()

# New data for rows 2-11
$newData = @(
    @(-0.0282525178045034, 0.0256563406437635, 0.06856962293386459),
    @(0.0575740486383438, -0.064446285367012, 0.1545489132404327),
    @(-0.0140499006956815, 0.0291688162833452, 0.0418442711234092),
    @(-0.0314595587551593, -0.0210748501121997, 0.0074830991216003),
    @(0.0100792767480015, -0.0314595587551593, -0.0143553335219621),
    @(-0.00137444678694, 0.0216857157647609, -0.016951510682702),
    @(0.0010690141934901, 0.0012217304902151, -0.0684169083833694),
    @(-0.0058032199740409, -0.0200058370828628, -0.07605272531509399),
    @(-0.0381790772080421, -0.0108428578823804, 0.0048869219608604),
    @(-0.0554360225796699, 0.0329867228865623, 0.012980886735022)
)

$r = 2
foreach ($row in $newData) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $r++
}
